# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across the resume body.
#
# Strategy: locate each target paragraph by its full (pre-edit) text via
# Find.Execute scoped to $d.Content, which yields an exact Start/End
# character range for that sentence. Then, working against the plain text
# captured from that same range, compute the character offsets of each
# metric substring and apply Bold + Color(#2C3E50) directly to the
# corresponding sub-range via $d.Range(start, end). Word automatically
# splits the run(s) as needed, matching the target OOXML (separate <w:r>
# runs with <w:rPr><w:b/><w:color w:val="2C3E50"/></w:rPr> around each
# metric, xml:space="preserve" added to the surrounding plain runs).

$d = $word.ActiveDocument

function RGBColor($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }
$highlightColor = RGBColor 0x2C 0x3E 0x50

function Highlight-Metrics($doc, $paragraphText, $segments) {
    $findRange = $doc.Content
    $found = $findRange.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $paragraphText"
        return
    }
    $rangeStart = $findRange.Start
    $plainText = $findRange.Text

    foreach ($seg in $segments) {
        $idx = $plainText.IndexOf($seg)
        if ($idx -lt 0) {
            Write-Output "SEGMENT NOT FOUND: '$seg' in '$paragraphText'"
            continue
        }
        $segStart = $rangeStart + $idx
        $segEnd = $segStart + $seg.Length
        $sub = $doc.Range($segStart, $segEnd)
        $sub.Font.Bold = 1
        $sub.Font.Color = $highlightColor
    }
}

Highlight-Metrics $d "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" @("23%", "64%")

Highlight-Metrics $d "• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes" @("±4.2%", "±2.1%", "71%", "87%")

Highlight-Metrics $d "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis" @("73.5%", "`$4.7M")

Highlight-Metrics $d "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion" @("`$2")

Highlight-Metrics $d "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%" @("57%")

Highlight-Metrics $d "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%" @("73.5%")

Highlight-Metrics $d "• `$4.7M savings enabled nonprofit access" @("`$4.7M")

Highlight-Metrics $d "• 178% accuracy improvement in racial classification algorithms" @("178%")
